# Update "想去人数" (want-to-go count, column F) figures across the 展览, 演出
# and 全部类型 sheets to reflect the latest scrape, per commit
# "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Cells.Item(2, 6).Value = 14
$wsExhibit.Cells.Item(5, 6).Value = 1007
$wsExhibit.Cells.Item(7, 6).Value = 2531
$wsExhibit.Cells.Item(9, 6).Value = 1236
$wsExhibit.Cells.Item(10, 6).Value = 905
$wsExhibit.Cells.Item(12, 6).Value = 908
$wsExhibit.Cells.Item(13, 6).Value = 1123
$wsExhibit.Cells.Item(17, 6).Value = 726
$wsExhibit.Cells.Item(18, 6).Value = 769
$wsExhibit.Cells.Item(19, 6).Value = 189
$wsExhibit.Cells.Item(20, 6).Value = 481
$wsExhibit.Cells.Item(21, 6).Value = 1108
$wsExhibit.Cells.Item(23, 6).Value = 593
$wsExhibit.Cells.Item(24, 6).Value = 586
$wsExhibit.Cells.Item(26, 6).Value = 298
$wsExhibit.Cells.Item(27, 6).Value = 296
$wsExhibit.Cells.Item(29, 6).Value = 351
$wsExhibit.Cells.Item(30, 6).Value = 4289
$wsExhibit.Cells.Item(36, 6).Value = 144
$wsExhibit.Cells.Item(37, 6).Value = 1596
$wsExhibit.Cells.Item(38, 6).Value = 438
$wsExhibit.Cells.Item(40, 6).Value = 83
$wsExhibit.Cells.Item(41, 6).Value = 140
$wsExhibit.Cells.Item(44, 6).Value = 127
$wsExhibit.Cells.Item(45, 6).Value = 132
$wsExhibit.Cells.Item(46, 6).Value = 96

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(8, 6).Value = 16
$wsShow.Cells.Item(13, 6).Value = 17

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(4, 6).Value = 14
$wsAll.Cells.Item(7, 6).Value = 1007
$wsAll.Cells.Item(8, 6).Value = 2531
$wsAll.Cells.Item(10, 6).Value = 1236
$wsAll.Cells.Item(11, 6).Value = 905
$wsAll.Cells.Item(13, 6).Value = 908
$wsAll.Cells.Item(14, 6).Value = 1123
$wsAll.Cells.Item(18, 6).Value = 726
$wsAll.Cells.Item(19, 6).Value = 769
$wsAll.Cells.Item(20, 6).Value = 189
$wsAll.Cells.Item(21, 6).Value = 481
$wsAll.Cells.Item(22, 6).Value = 1108
$wsAll.Cells.Item(25, 6).Value = 593
$wsAll.Cells.Item(26, 6).Value = 586
$wsAll.Cells.Item(28, 6).Value = 296
$wsAll.Cells.Item(30, 6).Value = 351
$wsAll.Cells.Item(31, 6).Value = 4289
$wsAll.Cells.Item(36, 6).Value = 144
$wsAll.Cells.Item(37, 6).Value = 1596
$wsAll.Cells.Item(38, 6).Value = 438
$wsAll.Cells.Item(39, 6).Value = 17
$wsAll.Cells.Item(40, 6).Value = 17
$wsAll.Cells.Item(42, 6).Value = 83
$wsAll.Cells.Item(43, 6).Value = 140
$wsAll.Cells.Item(46, 6).Value = 127
$wsAll.Cells.Item(48, 6).Value = 96
